# Add the new stimuli-detail content at the bottom of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Practice rows (2-5) gain a "pair_kind" value ("generic") in column J.
$ws.Range("J2").Value = "generic"
$ws.Range("J3").Value = "generic"
$ws.Range("J4").Value = "generic"
$ws.Range("J5").Value = "generic"

# New "stim details" block starting at row 27.
$ws.Range("A27").Value = "stim details"

# Header row for the new block.
$ws.Range("A28").Value = "month"
$ws.Range("B28").Value = "word_type"
$ws.Range("C28").Value = "need_audio"
$ws.Range("D28").Value = "need_image"
$ws.Range("E28").Value = "word"
$ws.Range("F28").Value = "count"
$ws.Range("G28").Value = "find images"

# Data rows 29-36: month number in A, word_type (video/audio) in B.
$months = @(6, 6, 7, 7, 6, 6, 7, 7)
$types  = @("video", "video", "video", "video", "audio", "audio", "audio", "audio")

for ($i = 0; $i -lt $months.Length; $i++) {
    $row = 29 + $i
    $ws.Range("A$row").Value = $months[$i]
    $ws.Range("B$row").Value = $types[$i]
}
